# Fix Training Data Issue
# The "Date" column (BF) values were recorded one day off because of how
# NBA stats were shown. Correct the raw "5-26-2012-13" label text to the
# proper ISO-style date text "2013-05-26" for every data row (BF2:BF31).
#
# NumberFormat is forced to Text ("@") first so Excel stores the literal
# string instead of silently reinterpreting the ISO-looking value as a
# date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("BF2:BF31")
$dataRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "2013-05-26"
}
